# Auto-generated script to apply numeric corrections to Ultima_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 990.36365
$ws.Range("I32").Value = 591
$ws.Range("J32").Value = 1030.3
$ws.Range("K32").Value = 591
$ws.Range("L32").Value = 1030.3
$ws.Range("M32").Value = -265
$ws.Range("N32").Value = -1682.3
$ws.Range("H33").Value = 33333694
$ws.Range("I33").Value = 34483120
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 34483120
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -34482891
$ws.Range("N33").Value = -758
$ws.Range("H41").Value = 1150.9333
$ws.Range("I41").Value = 1226.5714
$ws.Range("J41").Value = 92
$ws.Range("K41").Value = 1226.5714
$ws.Range("L41").Value = 92
$ws.Range("M41").Value = -786.5714
$ws.Range("N41").Value = -972
$ws.Range("H98").Value = 1979.9166
$ws.Range("I98").Value = 1528.7778
$ws.Range("J98").Value = 3333.3333
$ws.Range("K98").Value = 1528.7778
$ws.Range("L98").Value = 3333.3333
$ws.Range("M98").Value = -30.77780000000007
$ws.Range("N98").Value = -6329.3333
$ws.Range("H122").Value = 1979.9166
$ws.Range("I122").Value = 1528.7778
$ws.Range("J122").Value = 3333.3333
$ws.Range("K122").Value = 4586.3334
$ws.Range("L122").Value = 9999.999899999999
$ws.Range("M122").Value = -2136.3334
$ws.Range("N122").Value = -14899.9999
$ws.Range("H132").Value = 3499.9688
$ws.Range("I132").Value = 3258.0322
$ws.Range("K132").Value = 9774.096600000001
$ws.Range("M132").Value = -7244.096600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6041.427
$ws.Range("I32").Value = 6598.672
$ws.Range("J32").Value = 3552.4
$ws.Range("K32").Value = 6598.672
$ws.Range("L32").Value = 3552.4
$ws.Range("M32").Value = -6311.672
$ws.Range("N32").Value = -4126.4
$ws.Range("H74").Value = 14288478
$ws.Range("I74").Value = 19232034
$ws.Range("K74").Value = 19232034
$ws.Range("M74").Value = -19231160
$ws.Range("H77").Value = 14288478
$ws.Range("I77").Value = 19232034
$ws.Range("K77").Value = 96160170
$ws.Range("M77").Value = -96155802
$ws.Range("H122").Value = 7743.353
$ws.Range("I122").Value = 8966.143
$ws.Range("J122").Value = 2037
$ws.Range("K122").Value = 26898.429
$ws.Range("L122").Value = 6111
$ws.Range("M122").Value = -24448.429
$ws.Range("N122").Value = -11011

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 447.83334
$ws.Range("I22").Value = 466.66666
$ws.Range("J22").Value = 429
$ws.Range("K22").Value = 466.66666
$ws.Range("L22").Value = 429
$ws.Range("M22").Value = -293.66666
$ws.Range("N22").Value = -775

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4281
$ws.Range("I58").Value = 749.6667
$ws.Range("J58").Value = 6399.8
$ws.Range("K58").Value = 749.6667
$ws.Range("L58").Value = 6399.8
$ws.Range("M58").Value = -546.6667
$ws.Range("N58").Value = -6805.8
$ws.Range("H94").Value = 4474.0557
$ws.Range("I94").Value = 3060.5715
$ws.Range("J94").Value = 5373.5454
$ws.Range("K94").Value = 3060.5715
$ws.Range("L94").Value = 5373.5454
$ws.Range("M94").Value = -2609.5715
$ws.Range("N94").Value = -6275.5454
$ws.Range("H107").Value = 749.53845
$ws.Range("J107").Value = 1250
$ws.Range("L107").Value = 1250
$ws.Range("N107").Value = -5090
$ws.Range("H136").Value = 4281
$ws.Range("I136").Value = 749.6667
$ws.Range("J136").Value = 6399.8
$ws.Range("K136").Value = 2249.0001
$ws.Range("L136").Value = 19199.4
$ws.Range("M136").Value = 300.9998999999998
$ws.Range("N136").Value = -24299.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 89
$ws.Range("I8").Value = 89
$ws.Range("K8").Value = 267
$ws.Range("M8").Value = -128
$ws.Range("H97").Value = 9885.583000000001
$ws.Range("I97").Value = 34147.668
$ws.Range("J97").Value = 1798.2222
$ws.Range("K97").Value = 102443.004
$ws.Range("L97").Value = 5394.6666
$ws.Range("M97").Value = -101947.004
$ws.Range("N97").Value = -6386.6666
$ws.Range("H113").Value = 1209.6154
$ws.Range("I113").Value = 455.75
$ws.Range("J113").Value = 1855.7858
$ws.Range("K113").Value = 1367.25
$ws.Range("L113").Value = 5567.357400000001
$ws.Range("M113").Value = 802.75
$ws.Range("N113").Value = -9907.357400000001
$ws.Range("H121").Value = 977.0833
$ws.Range("J121").Value = 1566.4286
$ws.Range("L121").Value = 4699.2858
$ws.Range("N121").Value = -7319.2858
$ws.Range("H132").Value = 458
$ws.Range("I132").Value = 458
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4122
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1592
$ws.Range("N132").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6668737.5
$ws.Range("I122").Value = 8334672
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 25004016
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = -25001566
$ws.Range("N122").Value = -19897
$ws.Range("H126").Value = 4375.4165
$ws.Range("I126").Value = 3025.2
$ws.Range("J126").Value = 5339.857
$ws.Range("K126").Value = 9075.599999999999
$ws.Range("L126").Value = 16019.571
$ws.Range("M126").Value = -6605.599999999999
$ws.Range("N126").Value = -20959.571
$ws.Range("H132").Value = 3465.1333
$ws.Range("I132").Value = 2340.3438
$ws.Range("J132").Value = 6233.846
$ws.Range("K132").Value = 7021.0314
$ws.Range("L132").Value = 18701.538
$ws.Range("M132").Value = -4491.0314
$ws.Range("N132").Value = -23761.538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4599.7295
$ws.Range("I7").Value = 4493.5293
$ws.Range("J7").Value = 4690
$ws.Range("K7").Value = 4493.5293
$ws.Range("L7").Value = 4690
$ws.Range("M7").Value = -4381.5293
$ws.Range("N7").Value = -4914
$ws.Range("H40").Value = 4104.3794
$ws.Range("I40").Value = 4708
$ws.Range("J40").Value = 3361.4614
$ws.Range("K40").Value = 4708
$ws.Range("L40").Value = 3361.4614
$ws.Range("M40").Value = -4572
$ws.Range("N40").Value = -3633.4614
$ws.Range("H61").Value = 1349.5238
$ws.Range("I61").Value = 1213.9231
$ws.Range("J61").Value = 1569.875
$ws.Range("K61").Value = 1213.9231
$ws.Range("L61").Value = 1569.875
$ws.Range("M61").Value = -1011.9231
$ws.Range("N61").Value = -1973.875
$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25450
$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26560
$ws.Range("H113").Value = 1349.5238
$ws.Range("I113").Value = 1213.9231
$ws.Range("J113").Value = 1569.875
$ws.Range("K113").Value = 1213.9231
$ws.Range("L113").Value = 1569.875
$ws.Range("M113").Value = 956.0769
$ws.Range("N113").Value = -5909.875
$ws.Range("H126").Value = 4599.7295
$ws.Range("I126").Value = 4493.5293
$ws.Range("J126").Value = 4690
$ws.Range("K126").Value = 13480.5879
$ws.Range("L126").Value = 14070
$ws.Range("M126").Value = -11010.5879
$ws.Range("N126").Value = -19010
$ws.Range("H136").Value = 22735428
$ws.Range("I136").Value = 29412990
$ws.Range("J136").Value = 31720
$ws.Range("K136").Value = 88238970
$ws.Range("L136").Value = 95160
$ws.Range("M136").Value = -88236420
$ws.Range("N136").Value = -100260

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1519.8837
$ws.Range("I126").Value = 1176.7941
$ws.Range("K126").Value = 3530.3823
$ws.Range("M126").Value = -1060.3823
$ws.Range("H136").Value = 1238.6333
$ws.Range("I136").Value = 1160.5834
$ws.Range("K136").Value = 3481.7502
$ws.Range("M136").Value = -931.7501999999999
